# Updated symbol list on Sat Dec 31 06:26:58 UTC 2022 with GitHub Actions
# This script applies the scraped price/volume/link/coin-name refresh to
# the "cryptos" sheet, matching the upstream GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates -------------------------------------------
# NOTE: values are text in this sheet (t="inlineStr"), not numbers, so we
# prefix numeric-looking strings with a leading apostrophe to force Excel
# to store them as literal text and avoid floating point re-formatting.

$ws.Cells.Item(2,4).Value = "'245.08"  # D2
$ws.Cells.Item(4,4).Value = "'5.109"  # D4
$ws.Cells.Item(5,4).Value = "'0.05590"  # D5
$ws.Cells.Item(6,4).Value = "'6.474"  # D6
$ws.Cells.Item(7,4).Value = "'3.015"  # D7
$ws.Cells.Item(8,4).Value = "'0.8183"  # D8
$ws.Cells.Item(9,4).Value = "'0.8414"  # D9
$ws.Cells.Item(10,4).Value = "'0.1340"  # D10
$ws.Cells.Item(11,2).Value = "BitrueCoin"  # B11
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"  # C11
$ws.Cells.Item(11,4).Value = "'0.02856"  # D11
$ws.Cells.Item(11,5).Value = "10BitrueCoinBTR"  # E11
$ws.Cells.Item(12,2).Value = "BitMartToken"  # B12
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"  # C12
$ws.Cells.Item(12,4).Value = "'0.09373"  # D12
$ws.Cells.Item(12,5).Value = "11BitMartTokenBMX"  # E12
$ws.Cells.Item(13,2).Value = "BitForexToken"  # B13
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"  # C13
$ws.Cells.Item(13,4).Value = "'0.001516"  # D13
$ws.Cells.Item(13,5).Value = "12BitForexTokenBF"  # E13
$ws.Cells.Item(14,2).Value = "One"  # B14
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"  # C14
$ws.Cells.Item(14,4).Value = "'0.0005970"  # D14
$ws.Cells.Item(14,5).Value = "13OneONE"  # E14
$ws.Cells.Item(15,2).Value = "TigerCash"  # B15
$ws.Cells.Item(15,3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"  # C15
$ws.Cells.Item(15,4).Value = "'0.006224"  # D15
$ws.Cells.Item(15,5).Value = "14TigerCashTCH"  # E15
$ws.Cells.Item(16,2).Value = "LEO"  # B16
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"  # C16
$ws.Cells.Item(16,4).Value = "'3.523"  # D16
$ws.Cells.Item(16,5).Value = "15LEOLEO"  # E16
$ws.Cells.Item(17,2).Value = "BTSEToken"  # B17
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"  # C17
$ws.Cells.Item(17,4).Value = "'2.090"  # D17
$ws.Cells.Item(17,5).Value = "16BTSETokenBTSE"  # E17
$ws.Cells.Item(18,2).Value = "BitpandaEcosystemToken"  # B18
$ws.Cells.Item(18,3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"  # C18
$ws.Cells.Item(18,4).Value = "'0.3179"  # D18
$ws.Cells.Item(18,5).Value = "17BitpandaEcosystemTokenBEST"  # E18
$ws.Cells.Item(19,2).Value = "MandalaExchangeToken"  # B19
$ws.Cells.Item(19,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"  # C19
$ws.Cells.Item(19,4).Value = "'0.06954"  # D19
$ws.Cells.Item(19,5).Value = "18MandalaExchangeTokenMDX"  # E19
$ws.Cells.Item(20,2).Value = "LiechtensteinCryptoassetsExchange"  # B20
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"  # C20
$ws.Cells.Item(20,4).Value = "'0.03209"  # D20
$ws.Cells.Item(20,5).Value = "19LiechtensteinCryptoassetsExchangeLCX"  # E20
$ws.Cells.Item(22,4).Value = "'3.742"  # D22
$ws.Cells.Item(23,4).Value = "'0.04700"  # D23
$ws.Cells.Item(24,4).Value = "'0.1375"  # D24
$ws.Cells.Item(25,4).Value = "'0.001248"  # D25
$ws.Cells.Item(27,5).Value = "26NitroExNTX"  # E27
$ws.Cells.Item(40,4).Value = "'0.03661"  # D40
$ws.Cells.Item(41,4).Value = "'0.1364"  # D41
$ws.Cells.Item(41,5).Value = "40BKEXTokenBKKBestin24h"  # E41
$ws.Cells.Item(42,2).Value = "CEJI"  # B42
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"  # C42
$ws.Cells.Item(42,4).Value = "'0.002536"  # D42
$ws.Cells.Item(42,5).Value = "41CEJICEJI"  # E42
$ws.Cells.Item(43,2).Value = "KickToken"  # B43
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"  # C43
$ws.Cells.Item(43,4).Value = "'0.003374"  # D43
$ws.Cells.Item(43,5).Value = "42KickTokenKICK"  # E43
$ws.Cells.Item(44,4).Value = "'0.007521"  # D44
$ws.Cells.Item(45,4).Value = "'0.00005309"  # D45
$ws.Cells.Item(47,4).Value = "'0.1330"  # D47
$ws.Cells.Item(48,4).Value = "'0.002121"  # D48
$ws.Cells.Item(49,4).Value = "'0.00002100"  # D49
$ws.Cells.Item(50,4).Value = "'0.0002000"  # D50
